$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2201369519164543
$ws.Range("C2").Value = 0.9837354800394689
$ws.Range("D2").Value = 0.3575944958890494
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
$ws.Range("G2").Value = 0.1395347341502202
$ws.Range("H2").Value = 0.991
